# Update notification for order placing:
# Several backlog items (Product List, Cart Details, Cart Completion (Place
# Order) and Notification) have their Status moved from "Starting" to
# "Completed". Mirror the formatting already used by the other "Completed"
# rows on the BackLog sheet (e.g. row 33) instead of inventing a new style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BackLog")

$template = $ws.Range("H33")
$rows = @(11, 13, 18, 31)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 8)   # Column H = Status
    $cell.Value = "Completed"
    $template.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}

# Leave the selection where the author last worked, on one of the
# newly-updated Status cells.
$ws.Range("H18").Select() | Out-Null
